# feat(webapp): product filter work
#
# Rewrites the "schedule" sheet: shifts the weekly dates forward, replaces
# the TODO column text with the new set of tasks, turns the old "nyár" /
# "2. félév:" block into a red "Beadási határidő" callout, and replaces the
# trailing notes with the "Nyáron elkészült:" summary list. Rows 23-24 and
# the whole D column (the old "csúszás" markers) are dropped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Drop column D entirely (the "csúszás" markers on rows 6-8 and the
#    note on row 15) and the two trailing rows (23-24) that no longer
#    exist in the new schedule.
# ---------------------------------------------------------------------
$ws.Range("D1:D24").ClearContents()
$ws.Range("A23:C24").ClearContents()

# ---------------------------------------------------------------------
# 2. Weekly rows (2-15): keep the "# hét" labels in column A, push every
#    date forward by 203 days (29 weeks) and rewrite the TODO text.
# ---------------------------------------------------------------------
$dates = @(44809, 44816, 44823, 44830, 44837, 44844, 44851, 44858, 44865, 44872, 44879, 44886, 44893, 44900)
$todos = @(
  "Termék filterezés",
  "Termék filterezés, webes fizetés",
  "Webes fizetés, tesztek készítése",
  "Android kliens - architektura megtervezés, app skeleton létrehozása",
  "Android kliens - kezdőképernyő, autentikáció megvalósítás, lokális adatbázis megvalósítás",
  "Android kliens - termékek, kosár, checkout képernyők megvalósítása",
  "Android kliens - termék filterezés megvalósítása",
  "Diplomamunka írás",
  "Diplomamunka írás",
  "Diplomamunka írás",
  "Diplomamunka írás",
  "Diplomamunka írás",
  "Diplomamunka írás",
  "Diplomamunka írás"
)

for ($i = 0; $i -lt 14; $i++) {
  $row = 2 + $i
  $ws.Cells.Item($row, 2).Value = $dates[$i]
  $ws.Cells.Item($row, 3).Value = $todos[$i]
}

# ---------------------------------------------------------------------
# 3. Row 16 becomes the "Beadási határidő" / "2022.12.09, 12 óra" callout
#    in red, date-formatted text. Start from a clean (General) format so
#    the red-font style doesn't inherit A16's old date format, then build
#    the new format (red font, numFmtId 14 "mm-dd-yy") and copy that exact
#    format onto B16 so both cells share one style. Finally clear the old
#    C16 text.
# ---------------------------------------------------------------------
$a16 = $ws.Range("A16")
$a16.ClearFormats()
$a16.Value = "Beadási határidő"
$a16.Font.Color = 255
$a16.NumberFormat = "mm-dd-yy"

$a16.Copy()
$b16 = $ws.Range("B16")
$b16.PasteSpecial(-4122)
$b16.Value = "2022.12.09, 12 óra"
$excel.CutCopyMode = 0

$ws.Range("C16").ClearContents()

# ---------------------------------------------------------------------
# 4. Row 17 goes blank (keeps A17's existing date-style, loses its text
#    and the old C17 note).
# ---------------------------------------------------------------------
$ws.Range("A17").ClearContents()
$ws.Range("C17").ClearContents()

# ---------------------------------------------------------------------
# 5. Row 18: new "Nyáron elkészült:" label in column A (same style as the
#    other A-column labels, e.g. A2) plus the first summary bullet in C.
# ---------------------------------------------------------------------
$ws.Range("A2").Copy()
$a18 = $ws.Range("A18")
$a18.PasteSpecial(-4122)
$excel.CutCopyMode = 0
$a18.Value = "Nyáron elkészült:"

$ws.Range("C18").Value = "teljes admin felület (rendelések megtekintése, termékeken, kategóriákon CRUD műveletek)"

# ---------------------------------------------------------------------
# 6. Remaining summary bullets (rows 19-22). Row 20's column A goes
#    blank but keeps its existing style.
# ---------------------------------------------------------------------
$ws.Range("C19").Value = "színek, anyagok, minták kezelése"

$ws.Range("A20").ClearContents()
$ws.Range("C20").Value = "backend és fronted pagination rendelésekhez (admin) és termékekhez (user)"

$ws.Range("C21").Value = "rendelések és termékek query-zése (filterezés, rendezés)"

$ws.Range("C22").Value = "termékek és kategóriák képeinek kezelése - backenden külön thumbnail és original lementése (3rd party könyvtár segítségével)"

# ---------------------------------------------------------------------
# 7. Match the author's final selection.
# ---------------------------------------------------------------------
$ws.Range("C10").Select()

Write-Output "applied schedule update"
